$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "derivatives_dir" column (G) and its sample value are no longer
# needed, so remove the entire column. Everything to its right
# (modality0, modality0.input_source, modality1, modality1.input_source
# and their data validations) shifts left by one column automatically.
$ws.Columns("G").Delete()

# Restore a sensible selection on the now-shifted sheet.
$ws.Range("A1:A4").Select() | Out-Null
